$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update column C (Förändrad / last-changed date) for all existing data rows
#    (rows 2-307) from 45172 to 45175.
$ws.Range("C2:C307").Value = 45175

# 2. Row 307 gains an explicit row height (ht="15" customHeight="1").
$ws.Rows.Item(307).RowHeight = 15

# 3. Append a new data row (308) for case "A 41289-2023".
$ws.Range("A308").Value = "A 41289-2023"

$ws.Range("B308").Value = 45174
$ws.Range("B308").NumberFormat = "YYYY-MM-DD"

$ws.Range("C308").Value = 45175
$ws.Range("C308").NumberFormat = "YYYY-MM-DD"

$ws.Range("D308").Value = "HALLANDS LÄN"
$ws.Range("E308").Value = "LAHOLM"

$ws.Range("G308").Value = 0.4
$ws.Range("H308").Value = 0
$ws.Range("I308").Value = 0
$ws.Range("J308").Value = 0
$ws.Range("K308").Value = 0
$ws.Range("L308").Value = 0
$ws.Range("M308").Value = 0
$ws.Range("N308").Value = 0
$ws.Range("O308").Value = 0
$ws.Range("P308").Value = 0
$ws.Range("Q308").Value = 0

$ws.Range("R308").Value = ""
$ws.Range("R308").WrapText = $true
